$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92: H92, I92, K92, M92
$ws.Range("H92").Value = 1309.1875
$ws.Range("I92").Value = 1309.1875
$ws.Range("K92").Value = 1309.1875
$ws.Range("M92").Value = -61.1875

# Row 93: H93, J93, L93, N93
$ws.Range("H93").Value = 37020.723
$ws.Range("J93").Value = 37020.723
$ws.Range("L93").Value = 37020.723
$ws.Range("N93").Value = -42012.723

# Row 112: H112, J112, L112, N112
$ws.Range("H112").Value = 1356.75
$ws.Range("J112").Value = 1371.0256
$ws.Range("L112").Value = 4113.0768
$ws.Range("N112").Value = -6329.0768

# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Range("H113").Value = 6199.0713
$ws.Range("I113").Value = 1731.6666
$ws.Range("J113").Value = 7417.4546
$ws.Range("K113").Value = 1731.6666
$ws.Range("L113").Value = 7417.4546
$ws.Range("M113").Value = 1522.3334
$ws.Range("N113").Value = -13925.4546

# Row 129: H129, J129, L129, N129
$ws.Range("H129").Value = 1739.439
$ws.Range("J129").Value = 1874.973
$ws.Range("L129").Value = 5624.919
$ws.Range("N129").Value = -15624.919

# Row 137: H137, I137, J137, K137, L137, M137, N137
$ws.Range("H137").Value = 735724.9399999999
$ws.Range("I137").Value = 1908356.9
$ws.Range("J137").Value = 2829.9
$ws.Range("K137").Value = 5725070.699999999
$ws.Range("L137").Value = 8489.700000000001
$ws.Range("M137").Value = -5722520.699999999
$ws.Range("N137").Value = -13589.7

# Row 138: H138, J138, L138, N138
$ws.Range("H138").Value = 2448.05
$ws.Range("J138").Value = 3155.2273
$ws.Range("L138").Value = 9465.6819
$ws.Range("N138").Value = -19745.6819

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32, I32, K32, M32
$ws.Range("H32").Value = 3500.9507
$ws.Range("I32").Value = 3145.348
$ws.Range("K32").Value = 3145.348
$ws.Range("M32").Value = -2858.348

# Row 61: H61, I61, J61, K61, L61, M61, N61
$ws.Range("H61").Value = 1908.1177
$ws.Range("I61").Value = 1924.1428
$ws.Range("J61").Value = 1833.3334
$ws.Range("K61").Value = 1924.1428
$ws.Range("L61").Value = 1833.3334
$ws.Range("M61").Value = -1712.1428
$ws.Range("N61").Value = -2257.3334

# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 3275.8157
$ws.Range("I122").Value = 3015.5862
$ws.Range("J122").Value = 4114.3335
$ws.Range("K122").Value = 9046.758600000001
$ws.Range("L122").Value = 12343.0005
$ws.Range("M122").Value = -6596.758600000001
$ws.Range("N122").Value = -17243.0005

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 4700.6665
$ws.Range("I132").Value = 3247.625
$ws.Range("J132").Value = 7606.75
$ws.Range("K132").Value = 9742.875
$ws.Range("L132").Value = 22820.25
$ws.Range("M132").Value = -7212.875
$ws.Range("N132").Value = -27880.25

# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 1908.1177
$ws.Range("I136").Value = 1924.1428
$ws.Range("J136").Value = 1833.3334
$ws.Range("K136").Value = 5772.428400000001
$ws.Range("L136").Value = 5500.0002
$ws.Range("M136").Value = -3222.428400000001
$ws.Range("N136").Value = -10600.0002

# Row 137: H137, J137, L137, N137
$ws.Range("H137").Value = 39836
$ws.Range("J137").Value = 39836
$ws.Range("L137").Value = 39836
$ws.Range("N137").Value = -50036

$ws = $wb.Worksheets.Item("BSM")
# Row 21: H21, J21, L21, N21
$ws.Range("H21").Value = 28271
$ws.Range("J21").Value = 28271
$ws.Range("L21").Value = 28271
$ws.Range("N21").Value = -28743

# Row 95: H95, J95, L95, N95
$ws.Range("H95").Value = 33750
$ws.Range("J95").Value = 33750
$ws.Range("L95").Value = 33750
$ws.Range("N95").Value = -39242

# Row 99: H99, I99, J99, K99, L99, M99, N99
$ws.Range("H99").Value = 2842.25
$ws.Range("I99").Value = 1491.8889
$ws.Range("J99").Value = 6893.3335
$ws.Range("K99").Value = 1491.8889
$ws.Range("L99").Value = 6893.3335
$ws.Range("M99").Value = 6.111100000000079
$ws.Range("N99").Value = -9889.333500000001

# Row 103: H103, J103, L103, N103
$ws.Range("H103").Value = 34000
$ws.Range("J103").Value = 34000
$ws.Range("L103").Value = 34000
$ws.Range("N103").Value = -36344

# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value = 3214.0862
$ws.Range("I134").Value = 1144.6207
$ws.Range("J134").Value = 5283.552
$ws.Range("K134").Value = 3433.8621
$ws.Range("L134").Value = 15850.656
$ws.Range("M134").Value = -898.8620999999998
$ws.Range("N134").Value = -20920.656

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31, I31, J31, K31, L31, M31, N31
$ws.Range("H31").Value = 252591.53
$ws.Range("I31").Value = 541370.4
$ws.Range("J31").Value = 3644.276
$ws.Range("K31").Value = 541370.4
$ws.Range("L31").Value = 3644.276
$ws.Range("M31").Value = -541075.4
$ws.Range("N31").Value = -4234.276

# Row 34: H34, I34, J34, K34, L34, M34, N34
$ws.Range("H34").Value = 252591.53
$ws.Range("I34").Value = 541370.4
$ws.Range("J34").Value = 3644.276
$ws.Range("K34").Value = 541370.4
$ws.Range("L34").Value = 3644.276
$ws.Range("M34").Value = -541168.4
$ws.Range("N34").Value = -4048.276

# Row 58: H58, I58, J58, K58, L58, M58, N58
$ws.Range("H58").Value = 2990.3928
$ws.Range("I58").Value = 1555.9546
$ws.Range("J58").Value = 8250
$ws.Range("K58").Value = 1555.9546
$ws.Range("L58").Value = 8250
$ws.Range("M58").Value = -1352.9546
$ws.Range("N58").Value = -8656

# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 2990.3928
$ws.Range("I136").Value = 1555.9546
$ws.Range("J136").Value = 8250
$ws.Range("K136").Value = 4667.8638
$ws.Range("L136").Value = 24750
$ws.Range("M136").Value = -2117.8638
$ws.Range("N136").Value = -29850

$ws = $wb.Worksheets.Item("CUL")
# Row 14: H14, I14, K14, M14
$ws.Range("H14").Value = 60.27778
$ws.Range("I14").Value = 60.27778
$ws.Range("K14").Value = 180.83334
$ws.Range("M14").Value = -7.833339999999993

# Row 68: H68, I68, J68, K68, L68, M68, N68
$ws.Range("H68").Value = 2326.963
$ws.Range("I68").Value = 845.3333
$ws.Range("J68").Value = 2950.8071
$ws.Range("K68").Value = 2535.9999
$ws.Range("L68").Value = 8852.4213
$ws.Range("M68").Value = -1724.9999
$ws.Range("N68").Value = -10474.4213

# Row 71: H71, I71, J71, K71, L71, M71, N71
$ws.Range("H71").Value = 2326.963
$ws.Range("I71").Value = 845.3333
$ws.Range("J71").Value = 2950.8071
$ws.Range("K71").Value = 7607.9997
$ws.Range("L71").Value = 26557.2639
$ws.Range("M71").Value = -3551.9997
$ws.Range("N71").Value = -34669.2639

# Row 75: H75, I75, J75, K75, L75, M75, N75
$ws.Range("H75").Value = 2053.5
$ws.Range("I75").Value = 853.25
$ws.Range("J75").Value = 3253.75
$ws.Range("K75").Value = 2559.75
$ws.Range("L75").Value = 9761.25
$ws.Range("M75").Value = -1561.75
$ws.Range("N75").Value = -11757.25

# Row 78: H78, I78, J78, K78, L78, M78, N78
$ws.Range("H78").Value = 2053.5
$ws.Range("I78").Value = 853.25
$ws.Range("J78").Value = 3253.75
$ws.Range("K78").Value = 7679.25
$ws.Range("L78").Value = 29283.75
$ws.Range("M78").Value = -2687.25
$ws.Range("N78").Value = -39267.75

# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Range("H113").Value = 3572194.8
$ws.Range("I113").Value = 654.0952
$ws.Range("J113").Value = 8929506
$ws.Range("K113").Value = 1962.2856
$ws.Range("L113").Value = 26788518
$ws.Range("M113").Value = 207.7144000000001
$ws.Range("N113").Value = -26792858

# Row 131: H131, I131, J131, K131, L131, M131, N131
$ws.Range("H131").Value = 848.5
$ws.Range("I131").Value = 638
$ws.Range("J131").Value = 855.0103
$ws.Range("K131").Value = 1914
$ws.Range("L131").Value = 2565.0309
$ws.Range("M131").Value = 3126
$ws.Range("N131").Value = -12645.0309

$ws = $wb.Worksheets.Item("GSM")
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 4603.8335
$ws.Range("I122").Value = 2831
$ws.Range("J122").Value = 6376.6665
$ws.Range("K122").Value = 8493
$ws.Range("L122").Value = 19129.9995
$ws.Range("M122").Value = -6043
$ws.Range("N122").Value = -24029.9995

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 5813.857
$ws.Range("I132").Value = 3424.5
$ws.Range("J132").Value = 8999.666999999999
$ws.Range("K132").Value = 10273.5
$ws.Range("L132").Value = 26999.001
$ws.Range("M132").Value = -7743.5
$ws.Range("N132").Value = -32059.001

$ws = $wb.Worksheets.Item("LTW")
# Row 40: H40, I40, J40, K40, L40, M40, N40
$ws.Range("H40").Value = 7112
$ws.Range("I40").Value = 6501.3335
$ws.Range("J40").Value = 8333.333000000001
$ws.Range("K40").Value = 6501.3335
$ws.Range("L40").Value = 8333.333000000001
$ws.Range("M40").Value = -6365.3335
$ws.Range("N40").Value = -8605.333000000001

# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 6182.091
$ws.Range("I122").Value = 3274.5
$ws.Range("J122").Value = 7843.5713
$ws.Range("K122").Value = 9823.5
$ws.Range("L122").Value = 23530.7139
$ws.Range("M122").Value = -7373.5
$ws.Range("N122").Value = -28430.7139

# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 4643.08
$ws.Range("I132").Value = 3347.5
$ws.Range("J132").Value = 5839
$ws.Range("K132").Value = 10042.5
$ws.Range("L132").Value = 17517
$ws.Range("M132").Value = -7512.5
$ws.Range("N132").Value = -22577

# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 3090.3333
$ws.Range("I136").Value = 1304.95
$ws.Range("J136").Value = 4713.409
$ws.Range("K136").Value = 3914.85
$ws.Range("L136").Value = 14140.227
$ws.Range("M136").Value = -1364.85
$ws.Range("N136").Value = -19240.227

# Row 137: H137, J137, L137, N137
$ws.Range("H137").Value = 47791.668
$ws.Range("J137").Value = 47791.668
$ws.Range("L137").Value = 47791.668
$ws.Range("N137").Value = -57991.668

# Row 141: H141, J141, L141, N141
$ws.Range("H141").Value = 32080.223
$ws.Range("J141").Value = 32080.223
$ws.Range("L141").Value = 32080.223
$ws.Range("N141").Value = -42440.223

$ws = $wb.Worksheets.Item("WVR")
# Row 119: H119, J119, L119, N119
$ws.Range("H119").Value = 28000
$ws.Range("J119").Value = 28000
$ws.Range("L119").Value = 28000
$ws.Range("N119").Value = -37676

# Row 130: H130, J130, L130, N130
$ws.Range("H130").Value = 51179.668
$ws.Range("J130").Value = 51179.668
$ws.Range("L130").Value = 51179.668
$ws.Range("N130").Value = -61219.668

# Row 133: H133, J133, L133, N133
$ws.Range("H133").Value = 34425.082
$ws.Range("J133").Value = 34425.082
$ws.Range("L133").Value = 34425.082
$ws.Range("N133").Value = -44545.082

# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 3804.8928
$ws.Range("I136").Value = 1618
$ws.Range("K136").Value = 4854
$ws.Range("M136").Value = -2304
